$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-26 (columns A-I).
# Column A = model name (string), columns B-I = numeric metrics.
$data = @(
        @("model_9_9_12", 0.651583184210772, 0.7753551199508597, -0.08554840244577533, 0.4131451593922615, 0.3855947554111481, 0.341054767370224, 1.605809450149536, 0.9362332820892334),
        @("model_9_9_11", 0.6542959287863025, 0.7779888241899641, -0.06098909020762466, 0.4251904194489127, 0.3825925290584564, 0.3370562493801117, 1.569479823112488, 0.9170169830322266),
        @("model_9_9_7", 0.658824391158332, 0.7680263830246702, 0.07939183755550361, 0.4814250619736631, 0.377580851316452, 0.3521811962127686, 1.361819744110107, 0.8273035883903503),
        @("model_9_9_10", 0.6607129039750785, 0.7841320745110339, -0.03064284920125471, 0.4415238539239358, 0.3754908442497253, 0.3277296125888824, 1.524589896202087, 0.8909596800804138),
        @("model_9_9_9", 0.6626687118355219, 0.7842487333248664, -0.003201936494804158, 0.4535579753400545, 0.3733263313770294, 0.32755246758461, 1.483997464179993, 0.8717610836029053),
        @("model_9_9_13", 0.663085309795153, 0.8415495164439704, -0.2709395609568019, 0.365598976816951, 0.3728652596473694, 0.2405587285757065, 1.880051374435425, 1.012085676193237),
        @("model_9_9_8", 0.6652829663231834, 0.7874616534430574, 0.01910322110537732, 0.4649093551764861, 0.370433121919632, 0.3226746320724487, 1.451002359390259, 0.8536518216133118),
        @("model_9_9_5", 0.6677251379970661, 0.7090637832646504, 0.4616124294683552, 0.6184989891295667, 0.3677303493022919, 0.4416979551315308, 0.7964157462120056, 0.6086239814758301),
        @("model_9_9_0", 0.6721707901461442, 0.7489388316707017, 0.9486930182248843, 0.8511247096455985, 0.3628103733062744, 0.3811598122119904, 0.07589641213417053, 0.2375067472457886),
        @("model_9_9_4", 0.6730577171550403, 0.6919173484957147, 0.547550341933507, 0.6473589769224782, 0.3618287742137909, 0.4677295684814453, 0.6692911386489868, 0.5625824928283691),
        @("model_9_9_2", 0.6742678934658677, 0.7318535750021008, 0.6239641213317229, 0.700822261604479, 0.3604894280433655, 0.4070985615253448, 0.5562552213668823, 0.4772903919219971),
        @("model_9_9_3", 0.6745854856992193, 0.6814079480552728, 0.6279341268792257, 0.6771395789107342, 0.3601379692554474, 0.4836848974227905, 0.5503825545310974, 0.5150722861289978),
        @("model_9_9_6", 0.6785415976063043, 0.7361596081564205, 0.3921525974563483, 0.6018417382287786, 0.3557597100734711, 0.4005611836910248, 0.8991649150848389, 0.6351980566978455),
        @("model_9_9_16", 0.6822929288074036, 0.8351327352694201, -0.2414207410550644, 0.3752469760659487, 0.3516081273555756, 0.2503006458282471, 1.836385250091553, 0.9966937899589539),
        @("model_9_9_15", 0.6854933683355009, 0.8353212834940249, -0.2208993689181371, 0.3842958663077651, 0.3480661511421204, 0.2500144243240356, 1.806028842926025, 0.9822577238082886),
        @("model_9_9_14", 0.6864477840941351, 0.8345341686319201, -0.2069336168892359, 0.389991267962801, 0.3470099270343781, 0.2512094080448151, 1.785369873046875, 0.973171591758728),
        @("model_9_9_17", 0.6892793906036117, 0.8346157745766063, -0.2584654259115011, 0.3675495184575199, 0.3438761234283447, 0.2510855197906494, 1.86159884929657, 1.008973956108093),
        @("model_9_9_1", 0.6943035356393452, 0.7459214201546309, 0.9328226960389492, 0.842680125806698, 0.3383159041404724, 0.3857408165931702, 0.09937276691198349, 0.2509787678718567),
        @("model_9_9_18", 0.6967641178901627, 0.8304101345430084, -0.2414234812836484, 0.3728663529750784, 0.3355927765369415, 0.2574704885482788, 1.836389303207397, 1.000491619110107),
        @("model_9_9_19", 0.6999204153197562, 0.8338197321419742, -0.2498563410203429, 0.370905296835353, 0.3320996761322021, 0.2522940337657928, 1.848863840103149, 1.003620266914368),
        @("model_9_9_22", 0.700078710267364, 0.8264365864826092, -0.287131671867251, 0.3509191961364374, 0.3319244682788849, 0.2635031044483185, 1.904003739356995, 1.035504817962646),
        @("model_9_9_20", 0.7014166785922709, 0.834384636946696, -0.2527732728922691, 0.3699190966210506, 0.330443799495697, 0.2514364123344421, 1.853178739547729, 1.005193591117859),
        @("model_9_9_21", 0.7015614584324374, 0.8257127292138344, -0.2744063581436287, 0.3561094070511772, 0.3302834928035736, 0.2646020650863647, 1.88517951965332, 1.027224779129028),
        @("model_9_9_24", 0.7030945544433923, 0.8177275053695517, -0.2934687186702769, 0.3437695119953891, 0.3285868465900421, 0.276725172996521, 1.91337776184082, 1.046911001205444),
        @("model_9_9_23", 0.7036948668093669, 0.8181015562732558, -0.28417700065991, 0.3480093386485048, 0.3279224932193756, 0.2761573195457458, 1.899632930755615, 1.040147066116333)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
